$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("E3").Value = "  -4.44%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -10.37%  "
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("E17").Value = "  -4.02%  "
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("E22").Value = "  -5.55%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  +11.82%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("E29").Value = "  -8.75%  "
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  -5.41%  "
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  -4.20%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  -8.74%  "
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -8.15%  "
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("E49").Value = "  -11.41%  "
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("E51").Value = "  -3.19%  "
